$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 319.54166
$ws.Range("I6").Value2 = 290.5238
$ws.Range("J6").Value2 = 522.6667
$ws.Range("K6").Value2 = 871.5714
$ws.Range("L6").Value2 = 1568.0001
$ws.Range("M6").Value2 = -759.5714
$ws.Range("N6").Value2 = -1792.0001
$ws.Range("H8").Value2 = 398.875
$ws.Range("I8").Value2 = 298.7143
$ws.Range("K8").Value2 = 896.1428999999999
$ws.Range("M8").Value2 = -757.1428999999999
$ws.Range("H17").Value2 = 2497.3333
$ws.Range("J17").Value2 = 2497.3333
$ws.Range("L17").Value2 = 7491.999899999999
$ws.Range("N17").Value2 = -7827.999899999999
$ws.Range("H47").Value2 = 34899.25
$ws.Range("I47").Value2 = 33199
$ws.Range("J47").Value2 = 40000
$ws.Range("K47").Value2 = 33199
$ws.Range("L47").Value2 = 40000
$ws.Range("M47").Value2 = -32227
$ws.Range("N47").Value2 = -41944
$ws.Range("H112").Value2 = 839064.75
$ws.Range("I112").Value2 = 3196.5
$ws.Range("J112").Value2 = 991040.8
$ws.Range("K112").Value2 = 9589.5
$ws.Range("L112").Value2 = 2973122.4
$ws.Range("M112").Value2 = -8481.5
$ws.Range("N112").Value2 = -2975338.4
$ws.Range("H138").Value2 = 2323
$ws.Range("J138").Value2 = 2912.8333
$ws.Range("L138").Value2 = 8738.499899999999
$ws.Range("N138").Value2 = -19018.4999
$ws.Range("H141").Value2 = 56968.11
$ws.Range("I141").Value2 = 63714.188
$ws.Range("K141").Value2 = 191142.564
$ws.Range("M141").Value2 = -185962.564

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value2 = 1935.2174
$ws.Range("I97").Value2 = 1475.5
$ws.Range("J97").Value2 = 5000
$ws.Range("K97").Value2 = 1475.5
$ws.Range("L97").Value2 = 5000
$ws.Range("M97").Value2 = -979.5
$ws.Range("N97").Value2 = -5992
$ws.Range("H113").Value2 = 49398
$ws.Range("J113").Value2 = 49398
$ws.Range("L113").Value2 = 49398
$ws.Range("N113").Value2 = -58076
$ws.Range("H122").Value2 = 1250.7142
$ws.Range("I122").Value2 = 1088.8334
$ws.Range("J122").Value2 = 2222
$ws.Range("K122").Value2 = 3266.5002
$ws.Range("L122").Value2 = 6666
$ws.Range("M122").Value2 = -816.5001999999999
$ws.Range("N122").Value2 = -11566

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value2 = 1280.2727
$ws.Range("I5").Value2 = 2108.6
$ws.Range("J5").Value2 = 590
$ws.Range("K5").Value2 = 2108.6
$ws.Range("L5").Value2 = 590
$ws.Range("M5").Value2 = -1995.6
$ws.Range("N5").Value2 = -816
$ws.Range("H94").Value2 = 2416.9167
$ws.Range("I94").Value2 = 2865.6667
$ws.Range("K94").Value2 = 2865.6667
$ws.Range("M94").Value2 = -2414.6667
$ws.Range("H134").Value2 = 3473.4211
$ws.Range("I134").Value2 = 3473.4211
$ws.Range("K134").Value2 = 10420.2633
$ws.Range("M134").Value2 = -7885.263300000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 600
$ws.Range("I22").Value2 = 549.6667
$ws.Range("K22").Value2 = 549.6667
$ws.Range("M22").Value2 = -199.6667
$ws.Range("H43").Value2 = 21472.666
$ws.Range("J43").Value2 = 21472.666
$ws.Range("L43").Value2 = 21472.666
$ws.Range("N43").Value2 = -21840.666
$ws.Range("H58").Value2 = 2043.3846
$ws.Range("I58").Value2 = 2107.6
$ws.Range("J58").Value2 = 1829.3334
$ws.Range("K58").Value2 = 2107.6
$ws.Range("L58").Value2 = 1829.3334
$ws.Range("M58").Value2 = -1904.6
$ws.Range("N58").Value2 = -2235.3334
$ws.Range("H95").Value2 = 22561
$ws.Range("J95").Value2 = 22561
$ws.Range("L95").Value2 = 22561
$ws.Range("N95").Value2 = -28053
$ws.Range("H99").Value2 = 1622.6666
$ws.Range("I99").Value2 = 1444.4
$ws.Range("J99").Value2 = 2514
$ws.Range("K99").Value2 = 1444.4
$ws.Range("L99").Value2 = 2514
$ws.Range("M99").Value2 = 53.59999999999991
$ws.Range("N99").Value2 = -5510
$ws.Range("H101").Value2 = 21472.666
$ws.Range("J101").Value2 = 21472.666
$ws.Range("L101").Value2 = 21472.666
$ws.Range("N101").Value2 = -27962.666
$ws.Range("H105").Value2 = 3789.8572
$ws.Range("I105").Value2 = 2898
$ws.Range("J105").Value2 = 3938.5
$ws.Range("K105").Value2 = 2898
$ws.Range("L105").Value2 = 3938.5
$ws.Range("M105").Value2 = -1151
$ws.Range("N105").Value2 = -7432.5
$ws.Range("H126").Value2 = 1622.6666
$ws.Range("I126").Value2 = 1444.4
$ws.Range("J126").Value2 = 2514
$ws.Range("K126").Value2 = 4333.200000000001
$ws.Range("L126").Value2 = 7542
$ws.Range("M126").Value2 = -1863.200000000001
$ws.Range("N126").Value2 = -12482
$ws.Range("H132").Value2 = 1207.9286
$ws.Range("I132").Value2 = 1035.25
$ws.Range("K132").Value2 = 3105.75
$ws.Range("M132").Value2 = -575.75
$ws.Range("H136").Value2 = 2043.3846
$ws.Range("I136").Value2 = 2107.6
$ws.Range("J136").Value2 = 1829.3334
$ws.Range("K136").Value2 = 6322.799999999999
$ws.Range("L136").Value2 = 5488.0002
$ws.Range("M136").Value2 = -3772.799999999999
$ws.Range("N136").Value2 = -10588.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value2 = 187.42857
$ws.Range("I6").Value2 = 187.42857
$ws.Range("K6").Value2 = 562.28571
$ws.Range("M6").Value2 = -449.28571
$ws.Range("H10").Value2 = 411
$ws.Range("I10").Value2 = 366.5
$ws.Range("J10").Value2 = 500
$ws.Range("K10").Value2 = 1099.5
$ws.Range("L10").Value2 = 1500
$ws.Range("M10").Value2 = -960.5
$ws.Range("N10").Value2 = -1778
$ws.Range("H18").Value2 = 3738.4
$ws.Range("J18").Value2 = 2979.25
$ws.Range("L18").Value2 = 8937.75
$ws.Range("N18").Value2 = -9275.75
$ws.Range("H113").Value2 = 600.25
$ws.Range("I113").Value2 = 399.5
$ws.Range("J113").Value2 = 622.55554
$ws.Range("K113").Value2 = 1198.5
$ws.Range("L113").Value2 = 1867.66662
$ws.Range("M113").Value2 = 971.5
$ws.Range("N113").Value2 = -6207.66662
$ws.Range("H114").Value2 = 1347
$ws.Range("J114").Value2 = 4000
$ws.Range("L114").Value2 = 12000
$ws.Range("N114").Value2 = -18508
$ws.Range("H116").Value2 = 107393.734
$ws.Range("I116").Value2 = 149503.7
$ws.Range("K116").Value2 = 448511.1
$ws.Range("M116").Value2 = -445069.1

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("K4").Value2 = 0
$ws.Range("H18").Value2 = 15998.5
$ws.Range("J18").Value2 = 15998.5
$ws.Range("L18").Value2 = 15998.5
$ws.Range("N18").Value2 = -16584.5
$ws.Range("H80").Value2 = 2732
$ws.Range("J80").Value2 = 2737.25
$ws.Range("L80").Value2 = 2737.25
$ws.Range("N80").Value2 = -4733.25
$ws.Range("H83").Value2 = 2732
$ws.Range("J83").Value2 = 2737.25
$ws.Range("L83").Value2 = 13686.25
$ws.Range("N83").Value2 = -23670.25
$ws.Range("H97").Value2 = 1100.1428
$ws.Range("J97").Value2 = 1483.6666
$ws.Range("L97").Value2 = 1483.6666
$ws.Range("N97").Value2 = -2475.6666
$ws.Range("H113").Value2 = 3000
$ws.Range("J113").Value2 = 0
$ws.Range("L113").Value2 = 0
$ws.Range("H122").Value2 = 3050.238
$ws.Range("J122").Value2 = 2129.9
$ws.Range("L122").Value2 = 6389.700000000001
$ws.Range("N122").Value2 = -11289.7
$ws.Range("M4").ClearContents()
$ws.Range("N113").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 6166
$ws.Range("I40").Value2 = 4965.6665
$ws.Range("J40").Value2 = 7709.2856
$ws.Range("K40").Value2 = 4965.6665
$ws.Range("L40").Value2 = 7709.2856
$ws.Range("M40").Value2 = -4829.6665
$ws.Range("N40").Value2 = -7981.2856
$ws.Range("H61").Value2 = 2299.7273
$ws.Range("I61").Value2 = 1031.8889
$ws.Range("K61").Value2 = 1031.8889
$ws.Range("M61").Value2 = -829.8888999999999
$ws.Range("H82").Value2 = 3784.9285
$ws.Range("J82").Value2 = 3332
$ws.Range("L82").Value2 = 3332
$ws.Range("N82").Value2 = -4054
$ws.Range("H85").Value2 = 3784.9285
$ws.Range("J85").Value2 = 3332
$ws.Range("L85").Value2 = 3332
$ws.Range("N85").Value2 = -5828
$ws.Range("H103").Value2 = 23026
$ws.Range("J103").Value2 = 23026
$ws.Range("L103").Value2 = 23026
$ws.Range("N103").Value2 = -25370
$ws.Range("H113").Value2 = 2299.7273
$ws.Range("I113").Value2 = 1031.8889
$ws.Range("K113").Value2 = 1031.8889
$ws.Range("M113").Value2 = 1138.1111
$ws.Range("H141").Value2 = 100125
$ws.Range("J141").Value2 = 100000
$ws.Range("L141").Value2 = 100000
$ws.Range("N141").Value2 = -110360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value2 = 5334
$ws.Range("I3").Value2 = 5334
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 5334
$ws.Range("L3").Value2 = 0
$ws.Range("N3").Value2 = -5220
$ws.Range("H11").Value2 = 15244.5
$ws.Range("J11").Value2 = 25989
$ws.Range("L11").Value2 = 25989
$ws.Range("N11").Value2 = -26273
$ws.Range("H101").Value2 = 14013.25
$ws.Range("J101").Value2 = 14013.25
$ws.Range("L101").Value2 = 14013.25
$ws.Range("N101").Value2 = -20503.25
$ws.Range("H113").Value2 = 652.9091
$ws.Range("I113").Value2 = 585.4
$ws.Range("K113").Value2 = 1756.2
$ws.Range("M113").Value2 = 413.8000000000002
$ws.Range("H126").Value2 = 3657.7144
$ws.Range("I126").Value2 = 1868
$ws.Range("K126").Value2 = 5604
$ws.Range("M126").Value2 = -3134
$ws.Range("M3").ClearContents()
